# Apply the "rename extension" data fix to ProjectConfiguration.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the dataFile value (row 11, column B) from the old TestProject-specific
# file name to the generic example file name.
$ws.Range("B11").Value = "example_TimeValuesData.xlsx"

# Match the author's final cell selection in the saved workbook.
$ws.Range("B12").Select()
